# Add files via upload
# Rebuild the screen-requirements matrix: insert confirmation/graph screens,
# renumber downstream "N.xxx" references, extend the table, and trim the
# stray trailing blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Full row data (A..E) for rows 2..22 of the rebuilt matrix ---------
$rows = @(
    @(1,  "ログイン画面",             "ユーザーIDとパスワードを入力し、認証を行う", "-",                   "利用者・管理者共通"),
    @(2,  "ログアウト画面",           "システムからログアウトする",                 "各種メイン画面",       "-"),
    @(3,  "パスワードリセット画面",   "パスワードのリセットを行う",                 "1.ログイン画面",       "-"),
    @(4,  "利用者メイン画面",         "収支データや収支目標を管理するメニューを表示", "1.ログイン画面",      "-"),
    @(5,  "管理者メイン画面",         "ユーザー管理や収支目標設定支援を行うメニューを表示", "1.ログイン画面", "-"),
    @(6,  "ユーザー登録画面",         "新規ユーザーの情報を入力",                   "5.管理者メイン画面",   "管理者のみ"),
    @(7,  "ユーザー登録確認画面",     "入力内容を確認し、登録を確定する",           "6.ユーザー登録画面",   "-"),
    @(8,  "ユーザー編集画面",         "登録済みのユーザー情報を編集",               "9.ユーザー一覧画面",   "管理者のみ"),
    @(9,  "ユーザー一覧画面",         "登録済みのユーザーを一覧表示",               "5.管理者メイン画面",   "-"),
    @(10, "ユーザー削除確認画面",     "ユーザー削除の最終確認を行う",               "9.ユーザー一覧画面",   "管理者のみ"),
    @(11, "収支データ登録画面",       "収入・支出のデータを登録",                   "4.利用者メイン画面",   "-"),
    @(12, "収支データ登録確認画面",   "入力内容を確認し、登録を確定する",           "11.収支データ登録画面", "-"),
    @(13, "収支データ編集画面",       "登録済みの収支データを編集",                 "14.収支データ一覧画面", "-"),
    @(14, "収支データ一覧画面",       "収支データを一覧表示し、詳細確認・編集が可能", "4.利用者メイン画面",  "-"),
    @(15, "収支目標登録画面",         "収支目標を新規登録する",                     "4.利用者メイン画面",   "-"),
    @(16, "収支目標登録確認画面",     "入力内容を確認し、登録を確定する",           "15.収支目標登録画面",  "-"),
    @(17, "収支目標編集画面",         "登録済みの収支目標を編集",                   "18.収支目標一覧画面",  "-"),
    @(18, "収支目標一覧画面",         "設定済みの収支目標を一覧表示",               "4.利用者メイン画面",   "-"),
    @(19, "収支目標詳細画面",         "収支目標の達成状況を詳細表示",               "18.収支目標一覧画面",  "-"),
    @(21, "収支データグラフ画面",     "収支データをグラフで可視化",                 "4.利用者メイン画面",   "-"),
    @(22, "収支目標グラフ画面",       "収支目標の達成率をグラフで表示",             "4.利用者メイン画面",   "-")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $r = $r + 1
}

# --- Grow the table (表_1 -> Table_1) to cover the new rows -------------
$lo = $ws.ListObjects.Item(1)
$lo.Name = "Table_1"
$lo.Resize($ws.Range("A1:E22"))

# --- Carry the existing row formatting down onto the newly added rows ---
# (reuses the same style records instead of minting new ones)
$ws.Range("A2:E2").Copy()
$ws.Range("A17:E22").PasteSpecial(-4122)

# --- Drop the now-stray trailing blank row (sheet had 999 rows, now 998) ---
$ws.Rows.Item(999).Delete()
